$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Producto"
$ws.Range("D1").Value = "Disponibilidad"

$ws.Range("D4").Select()

$ws.Columns.Item(4).ColumnWidth = 12.125
